$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Narrow column I (Комментарий) from 50 to 30 characters wide. The COM
# layer snaps ColumnWidth to whole pixels before it round-trips back to
# the stored OOXML width, so 29.15 (which snaps to the same pixel width
# as a true 30) is used to land exactly on width="30" in the saved file.
$ws.Columns.Item(9).ColumnWidth = 29.15

# Append a new data row (row 2) under the header row.
$ws.Range("A2").Value = "1 новость"
$ws.Range("B2").Value = "Краснова Ксения Максимовна"
$ws.Range("C2").Value = "sdfsdfsdf"
$ws.Range("D2").Value = "Технология"
$ws.Range("E2").Value = "1а"
$ws.Range("F2").Value = "Спорт"
$ws.Range("G2").Value = "Школьный"
$ws.Range("H2").Value = "Лауреат"
# I2 (Комментарий) is left blank for this row.

# J2/K2/L2 hold date-like text that must stay literal strings, not get
# reinterpreted as serial dates, so force text formatting for the write
# and then drop the formatting override afterwards.
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "12.02.2007"
$ws.Range("J2").ClearFormats()

$ws.Range("K2").Value = "29.03.2023"
$ws.Range("L2").Value = "29.03.2023"
